# Started Program Module from my End.
#
# Adds a new "Program" worksheet (mirroring the layout used by the other
# module sheets, e.g. "Login") and wires it up from the "Export Summary"
# index sheet, the same way every other module sheet is referenced there.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Program" worksheet at the end of the workbook.
#    Copying the existing "Login" sheet gives the new sheet the same
#    5-column / 10-row layout, frozen-look formatting, page setup and
#    footer that every module sheet in this workbook shares.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$program = $wb.Worksheets.Item($wb.Worksheets.Count)
$program.Name = "Program"

# Freeze the header row/column, matching the top-left frozen pane used on
# this sheet.
$program.Activate()
$null = $program.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Header row.
$program.Range("A1").Value = "programDescription"
$program.Range("B1").Value = "programStatus"
$program.Range("C1").Value = "programId"
$program.Range("D1").ClearContents()

# Data row.
$program.Range("A2").Value = "Selenium Classes"
$program.Range("B2").Value = "Active"
# The copied cell previously held text, which would otherwise coerce this
# into a string "0" instead of the numeric 0 the source data uses.
$program.Range("C2").NumberFormat = "General"
$program.Range("C2").Value = 0
$program.Range("D2").ClearContents()

# ---------------------------------------------------------------------
# 2. Reference the new sheet from the "Export Summary" sheet, copying
#    the existing two-row block pattern used for every other module
#    (label row + "Table 1" / hyperlink row).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Export Summary")

$lastBlock = $summary.Range("B15:D16")
$newBlock = $summary.Range("B17:D18")

# Add the hyperlink to the target cell first so the subsequent Copy
# (which carries over the correct fonts/fills/borders from the previous
# block) is what determines the final cell style.
$summary.Hyperlinks.Add($summary.Range("D18"), "", "'Program'!R1C1", "", "Program") | Out-Null

$lastBlock.Copy($newBlock)

$summary.Range("B17").Value = "Program"
$summary.Range("D18").Value = "Program"
